$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values would otherwise be
# misinterpreted as numbers (losing significant trailing zeros / exact literal form).
foreach ($addr in @("D4", "D6", "D19", "D30", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.479.47"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.837.72"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "261.89"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.5385"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").Value = "0.3016"
$ws.Range("E8").Value = "  -6.86%  "
$ws.Range("D9").Value = "0.06896"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").Value = "17.65"
$ws.Range("E10").Value = "  -6.72%  "
$ws.Range("D11").Value = "0.7385"
$ws.Range("E11").Value = "  -5.54%  "
$ws.Range("D12").Value = "1.847.09"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "0.07139"
$ws.Range("E13").Value = "  -8.12%  "
$ws.Range("D14").Value = "89.16"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "4.995"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "0.000007900"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "26.507.07"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "2.078.47"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "4.593"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "5.985"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").Value = "142.81"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "2.183"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "1.718"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").Value = "17.02"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "111.19"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "4.250"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "0.08854"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").Value = "4.052"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("D33").Value = "0.04842"
$ws.Range("D34").Value = "2.928"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "0.7302"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "3.097"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "2.263"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").Value = "0.01719"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "0.4714"
$ws.Range("D41").Value = "0.9046"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "108.22"
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").Value = "5.901"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "7.404"
$ws.Range("E45").Value = "  -3.34%  "
$ws.Range("D46").Value = "0.1253"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "9.006"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "0.4079"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").Value = "34.82"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("D50").Value = "0.8940"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "0.05766"
$ws.Range("E51").Value = "  -2.16%  "
